$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("OWQ9DY4JY4", "646", "8.4.2015 г.", "IrregularIncome", "дадад"),
    @("77ERLFSCXO", "6456", "8.4.2015 г.", "RegularIncome", "хжхгхг"),
    @("3F8T2TJK6A", "-656", "8.4.2015 г.", "RegularIncome", "gfgfgfgfg"),
    @("9112TX7OJP", "-5454", "8.4.2015 г.", "RegularIncome", "cggfgfg"),
    @("Z9RCKW82NQ", "654565", "8.4.2015 г.", "IrregularIncome", "gfdgfgf"),
    @("CBZQY6LA92", "1000", "8.4.2015 г.", "IrregularIncome", "gfgfgf"),
    @("A9UY5NENEM", "55555", "8.4.2015 г.", "IrregularExpense", "tttttttt"),
    @("LP3949ZN78", "5656", "8.4.2015 г.", "RegularIncome", ""),
    @("8YH22NYVXY", "64646", "8.4.2015 г.", "RegularIncome", ""),
    @("PUXJPLW6XD", "6565", "8.4.2015 г.", "RegularIncome", ""),
    @("L1XVPLQD2R", "965", "8.4.2015 г.", "RegularIncome", ""),
    @("5I3AUFJ1EF", "5754", "8.4.2015 г.", "IrregularIncome", "hfhh"),
    @("1HI3P6ZLEM", "5464", "8.4.2015 г.", "IrregularExpense", "хфгхг")
)

$startRow = 13
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = "'" + $data[1]
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    if ($data[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = $data[4]
    }
}

Write-Output "Rows added"